$wb = $excel.ActiveWorkbook

# --- 1. Sheet "SearchBarData": A1 "Adidas" -> "ADIDAS" ---
$ws1 = $wb.Worksheets.Item("SearchBarData")
$ws1.Range("A1").Value = "ADIDAS"

# --- 2. Sheet "AssertData": add D1 (empty-string cell), extends used range to A1:D7 ---
$ws2 = $wb.Worksheets.Item("AssertData")
$ws2.Range("D1").Formula = '=""'

# --- 3. Sheet "CheckBoxData" is unchanged ---

# --- 4. Add new sheet "PageURLS" after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsUrls = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$wsUrls.Name = "PageURLS"

$wsUrls.Range("A1").Value = "https://sportsjam.in/"
$wsUrls.Range("A2").Value = "https://sportsjam.in/shopbrand"
$wsUrls.Range("A3").Value = "https://sportsjam.in/shopbrand/puma"
$wsUrls.Range("A4").Value = "https://sportsjam.in/ox_quickview/catalog_product/view/id/30663/"

$wsUrls.Range("C1").Value = "https://sportsjam.in/sports/cricket-equipment-store-online-india/cricket-shoes-shop-online-india"
$wsUrls.Range("C2").Value = "https://sportsjam.in/"
$wsUrls.Range("C3").Value = "https://sportsjam.in/badminton-shoes-online-india"
$wsUrls.Range("C4").Value = "https://sportsjam.in/sports/basketball-gear-online-india/basketball-shoes-online-india"
$wsUrls.Range("C5").Value = "https://sportsjam.in/sports/buy-sports-football/football-shoes-online-india"
$wsUrls.Range("C6").Value = "https://sportsjam.in/buy-running/buy-running-shoes-online-india"
$wsUrls.Range("C7").Value = "https://sportsjam.in/indoor-court-squash-shoes-online-india"

# --- 5. Add new sheet "DataToBeSent" after PageURLS ---
$wsSent = $wb.Worksheets.Add([Type]::Missing, $wsUrls)
$wsSent.Name = "DataToBeSent"

$wsSent.Range("A1").Value = "Puma"
$wsSent.Range("A2").Value = "test134@gmail.com"
$wsSent.Range("A3").Value = "test140@ga"

# --- 6. Add new sheet "AssertContainsData" after DataToBeSent ---
$wsAssertContains = $wb.Worksheets.Add([Type]::Missing, $wsSent)
$wsAssertContains.Name = "AssertContainsData"

$wsAssertContains.Range("A1").Value = "Puma"
$wsAssertContains.Range("A2").Value = "accessories"
$wsAssertContains.Range("A3").Value = "you must login or register to add items to your wishlist."
$wsAssertContains.Range("A4").Value = "THE ACCOUNT SIGN-IN WAS INCORRECT OR YOUR ACCOUNT IS DISABLED TEMPORARILY. PLEASE WAIT AND TRY AGAIN LATER."

Write-Host "Sheets:"
$wb.Worksheets | ForEach-Object { Write-Host $_.Name }
